# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Toscana, Especial / Primera) right
# after row 199, shifting the existing rows 200-224 down to 202-226.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 200.
$ws.Range("A200:A201").EntireRow.Insert()

# New row 200: Toscana / Especial
$ws.Range("A200").Value = 11
$ws.Range("B200").Value = "Vega Monumental Concepción"
$ws.Range("C200").Value = "Bíobío"
$ws.Range("D200").Value = 44918
$ws.Range("E200").Value = 8
$ws.Range("F200").Value = "Fruta"
$ws.Range("G200").Value = 100103
$ws.Range("H200").Value = "Frutos de hueso (carozo)"
$ws.Range("I200").Value = 100103004
$ws.Range("J200").Value = "Durazno"
$ws.Range("K200").Value = "Toscana"
$ws.Range("L200").Value = "Especial"
$ws.Range("M200").Value = 50
$ws.Range("N200").Value = 18000
$ws.Range("O200").Value = 18000
$ws.Range("P200").Value = 18000
$ws.Range("Q200").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R200").Value = "Región de O'Higgins"
$ws.Range("S200").Value = 1200
$ws.Range("T200").Value = 15

# New row 201: Toscana / Primera
$ws.Range("A201").Value = 11
$ws.Range("B201").Value = "Vega Monumental Concepción"
$ws.Range("C201").Value = "Bíobío"
$ws.Range("D201").Value = 44918
$ws.Range("E201").Value = 8
$ws.Range("F201").Value = "Fruta"
$ws.Range("G201").Value = 100103
$ws.Range("H201").Value = "Frutos de hueso (carozo)"
$ws.Range("I201").Value = 100103004
$ws.Range("J201").Value = "Durazno"
$ws.Range("K201").Value = "Toscana"
$ws.Range("L201").Value = "Primera"
$ws.Range("M201").Value = 50
$ws.Range("N201").Value = 16000
$ws.Range("O201").Value = 16000
$ws.Range("P201").Value = 16000
$ws.Range("Q201").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R201").Value = "Región de O'Higgins"
$ws.Range("S201").Value = 1067
$ws.Range("T201").Value = 15
